$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.914.14"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.25%  "
$ws.Range("D3").Value = "'2.591.41"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.76%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'520.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.84%  "
$ws.Range("D6").Value = "'140.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.37%  "
$ws.Range("D7").Value = "'0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").Value = "'0.565"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.30%  "
$ws.Range("D9").Value = "'2.599.06"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.78%  "
$ws.Range("D10").Value = "'6.50"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.10%  "
$ws.Range("E11").Value = "  +1.91%  "
$ws.Range("D12").Value = "'0.332"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.90%  "
$ws.Range("E13").Value = "  +2.56%  "
$ws.Range("D14").Value = "'3.047.46"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.66%  "
$ws.Range("D15").Value = "'58.870.74"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.12%  "
$ws.Range("E16").Value = "  +2.44%  "
$ws.Range("D17").Value = "'2.603.33"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.83%  "
$ws.Range("E18").Value = "  +0.72%  "
$ws.Range("D19").Value = "'339.05"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.47%  "
$ws.Range("E20").Value = "  +1.78%  "
$ws.Range("D21").Value = "'10.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.32%  "
$ws.Range("D22").Value = "'6.54"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.19%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").Value = "'66.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.97%  "
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("E26").Value = "  +1.68%  "
$ws.Range("D27").Value = "'0.997"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("D28").Value = "'7.10"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.05%  "
$ws.Range("D29").Value = "'0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("D30").Value = "'0.0₃0726"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.49%  "
$ws.Range("D31").Value = "'5.96"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.60%  "
$ws.Range("D32").Value = "'18.78"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.69%  "
$ws.Range("D33").Value = "'1.57"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("D34").Value = "'148.88"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.33%  "
$ws.Range("D35").Value = "'4.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.72%  "
$ws.Range("E36").Value = "  -0.11%  "
$ws.Range("D37").Value = "'36.28"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.95%  "
$ws.Range("D38").Value = "'0.837"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.92%  "
$ws.Range("D39").Value = "'1.46"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.14%  "
$ws.Range("D40").Value = "'0.825"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.53%  "
$ws.Range("D41").Value = "'3.52"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.59%  "
$ws.Range("D42").Value = "'0.997"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.17%  "
$ws.Range("D43").Value = "'276.09"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.34%  "
$ws.Range("D44").Value = "'10.74"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.94%  "
$ws.Range("D45").Value = "'0.0951"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.33%  "
$ws.Range("D46").Value = "'0.588"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.66%  "
$ws.Range("D47").Value = "'0.0523"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.70%  "
$ws.Range("D48").Value = "'18.60"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.37%  "
$ws.Range("D49").Value = "'1.985.83"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.56%  "
$ws.Range("D50").Value = "'0.0221"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.94%  "
$ws.Range("D51").Value = "'4.48"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.56%  "
